$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly observed Jira issue timestamps / values for several rows
$ws.Range("B2").Value = "2026-01-15T16:00:24.012-0500"
$ws.Range("D2").Value = "2026-01-15T16:05:42.215-0500"
$ws.Range("B3").Value = "2026-01-15T15:45:24.794-0500"
$ws.Range("D3").Value = "2026-01-15T16:31:21.086-0500"
$ws.Range("B4").Value = "2026-01-15T15:15:53.448-0500"
$ws.Range("D4").Value = "2026-01-15T15:20:41.719-0500"
$ws.Range("B5").Value = "2026-01-15T15:35:47.674-0500"
$ws.Range("C5").Value = "2026-01-15T15:47:01.091-0500"
$ws.Range("B6").Value = "2026-01-15T15:29:46.925-0500"
$ws.Range("D6").Value = "2026-01-15T16:30:08.186-0500"
$ws.Range("B7").Value = "2026-01-15T12:34:53.684-0500"
$ws.Range("C7").Value = "2026-01-15T13:22:11.843-0500"
$ws.Range("B8").Value = "2026-01-15T12:58:46.600-0500"
$ws.Range("B9").Value = "2026-01-15T19:51:49.848-0500"
$ws.Range("B10").Value = "2026-01-15T10:19:41.368-0500"
$ws.Range("C10").Value = "2026-01-15T11:55:35.561-0500"
$ws.Range("B11").Value = "2026-01-15T10:11:15.402-0500"
$ws.Range("C11").Value = "2026-01-15T11:55:00.479-0500"
$ws.Range("B12").Value = "2026-01-15T09:47:10.662-0500"
$ws.Range("C12").Value = "2026-01-15T10:04:52.006-0500"
$ws.Range("B13").Value = "2026-01-15T11:38:21.744-0500"
$ws.Range("D13").Value = "2026-01-15T15:56:11.732-0500"
$ws.Range("B14").Value = "2026-01-15T18:49:35.472-0500"
$ws.Range("B15").Value = "2026-01-14T23:57:30.381-0500"
$ws.Range("C15").Value = "2026-01-15T00:32:42.229-0500"
$ws.Range("B16").Value = "2026-01-14T23:42:10.594-0500"
$ws.Range("C16").Value = "2026-01-15T01:08:03.854-0500"
$ws.Range("B17").Value = "2026-01-14T17:29:36.581-0500"
$ws.Range("C17").Value = "2026-01-14T20:46:57.091-0500"
$ws.Range("B18").Value = "2026-01-14T15:29:24.568-0500"
$ws.Range("C18").Value = "2026-01-14T19:09:10.304-0500"
$ws.Range("B19").Value = "2026-01-14T15:28:01.055-0500"
$ws.Range("C19").Value = "2026-01-14T19:06:00.051-0500"
$ws.Range("D19").Value = "2026-01-14T19:08:01.554-0500"
$ws.Range("B22").Value = "2026-01-14T12:23:48.584-0500"
$ws.Range("C22").Value = "2026-01-14T20:47:23.431-0500"
$ws.Range("B54").Value = "2026-01-11T15:12:33.149-0500"
$ws.Range("C54").Value = "2026-01-11T21:01:05.151-0500"
$ws.Range("D54").Value = "2026-01-11T21:05:53.010-0500"
$ws.Range("B83").Value = "2026-01-07T13:26:07.581-0500"
$ws.Range("C83").Value = "2026-01-07T14:59:10.479-0500"
$ws.Range("D85").Value = "2026-01-15T11:09:50.333-0500"
$ws.Range("B107").Value = "2026-01-01T08:31:30.704-0500"
$ws.Range("C107").Value = "2026-01-01T09:04:50.071-0500"
$ws.Range("B120").Value = "2025-12-31T06:10:53.768-0500"
$ws.Range("C120").Value = "2025-12-31T10:38:05.820-0500"
$ws.Range("B125").Value = "2025-12-30T01:27:32.138-0500"
$ws.Range("C125").Value = "2025-12-30T03:12:53.621-0500"
$ws.Range("D125").Value = "2026-01-14T10:24:45.273-0500"
$ws.Range("B129").Value = "2025-12-31T04:03:15.646-0500"
$ws.Range("C129").Value = "2025-12-31T10:37:56.064-0500"
$ws.Range("B131").Value = "2025-12-29T08:49:51.016-0500"
$ws.Range("C131").Value = "2025-12-29T09:34:44.539-0500"
$ws.Range("B133").Value = "2025-12-27T07:43:51.377-0500"
$ws.Range("C133").Value = "2025-12-27T07:53:22.239-0500"
$ws.Range("D133").Value = "2026-01-14T08:57:05.462-0500"
$ws.Range("B136").Value = "2025-12-29T07:07:36.908-0500"
$ws.Range("C136").Value = "2025-12-29T08:41:38.437-0500"
$ws.Range("D136").Value = "2025-12-26T19:46:09.269-0500"
$ws.Range("D138").Value = "2026-01-14T22:11:26.800-0500"
$ws.Range("D157").Value = "2026-01-14T22:08:27.542-0500"

# These cells hold numeric-looking text (e.g. "190.09") that must remain text,
# not be auto-converted to a number. Force text format, set the value, then
# clear the formatting override so no extra style is left applied to the cell.
$ws.Range("H85").NumberFormat = "@"
$ws.Range("H85").Value = "190.09"
$ws.Range("H85").ClearFormats()
$ws.Range("H138").NumberFormat = "@"
$ws.Range("H138").Value = "497.00"
$ws.Range("H138").ClearFormats()
$ws.Range("H157").NumberFormat = "@"
$ws.Range("H157").Value = "572.31"
$ws.Range("H157").ClearFormats()

# Remove the 15 oldest rows (TPGSOC-1314437 .. TPGSOC-1312078) which have
# aged out of the tracked window; this shifts everything below up and
# updates the sheet dimension accordingly (A1:H193 -> A1:H178).
$ws.Range("A179:H193").EntireRow.Delete()
